# lectura excel y validación de proyectos repetidos
#
# Duplicate row 18 (AR044321 / Construction / CON / a building or road)
# into the next three rows (27, 28, 29) so the sheet holds repeated
# project entries, then leave the selection like Excel would after
# pasting the last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sourceRow = 18
$startRow  = 27
$endRow    = 29

$colA = $ws.Range("A$sourceRow").Value2
$colB = $ws.Range("B$sourceRow").Value2
$colC = $ws.Range("C$sourceRow").Value2
$colD = $ws.Range("D$sourceRow").Value2

for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Range("A$r").Value = $colA
    $ws.Range("B$r").Value = $colB
    $ws.Range("C$r").Value = $colC
    $ws.Range("D$r").Value = $colD
}

$ws.Application.Goto($ws.Range("A16"))
$ws.Range("A$endRow`:XFD$endRow").Select()
